$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old "Fiscal Year From" row (old row 3) to make
# room for the new "Account" and "Partner" filter fields. This pushes the
# existing rows (old 3..10) down to (new 5..12).
$ws.Range("A3:A4").EntireRow.Insert()

# New "Account" filter row (new row 3) - copy formatting from the row that
# used to be the "Fiscal Year From" row (now at row 5) so the label/input
# styling matches the rest of the filter block.
$ws.Range("A5:B5").Copy()
$ws.Range("A3:B3").PasteSpecial(-4122)
$ws.Range("A3").Value = "Account"
$ws.Range("B3").Value = ""

# New "Partner" filter row (new row 4)
$ws.Range("A5:B5").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Partner"
$ws.Range("B4").Value = ""

# Row 1 becomes a new blank spacer row styled like the title row.
$ws.Range("A2").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Value = ""

# Fix "Activity" typos in the header row (now row 12 after the insert).
$ws.Range("E12").Value = "Activity Group Name"
$ws.Range("F12").Value = "Activity Code"

# Match the upstream row-height tweak for the header row.
$ws.Rows.Item(12).RowHeight = 28.35

# Match the upstream selection (active cell moved to E2).
$ws.Range("E2").Select()
